$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = "Testing 2"
$ws.Range("B3").Value = "asdfdsf"
$ws.Range("C3").Value = "ABC"
$ws.Range("D3").Value = 250

$ws.Range("A4").Value = "Mike"
$ws.Range("B4").Value = "Some description"
$ws.Range("C4").Value = "USD"
$ws.Range("D4").Value = 400

$ws.Range("D5").Select()
